$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 22).Value = 1.17

# Row 3
$ws.Cells.Item(3, 6).Value = 2.5
$ws.Cells.Item(3, 7).Value = 2.72
$ws.Cells.Item(3, 9).Value = 3.7
$ws.Cells.Item(3, 12).Value = 1.49
$ws.Cells.Item(3, 13).Value = 1.1
$ws.Cells.Item(3, 16).Value = 1.53
$ws.Cells.Item(3, 17).Value = 2.68
$ws.Cells.Item(3, 19).Value = 4.8
$ws.Cells.Item(3, 23).Value = 1.58

# Row 5
$ws.Cells.Item(5, 6).Value = 1.2
$ws.Cells.Item(5, 7).Value = 1.21
$ws.Cells.Item(5, 8).Value = 13
$ws.Cells.Item(5, 10).Value = 8
$ws.Cells.Item(5, 11).Value = 9.800000000000001
$ws.Cells.Item(5, 15).Value = 1.13
$ws.Cells.Item(5, 16).Value = 2.98
$ws.Cells.Item(5, 23).Value = 5.7
$ws.Cells.Item(5, 29).Value = 22
$ws.Cells.Item(5, 33).Value = 14.5
$ws.Cells.Item(5, 37).Value = 16.5
$ws.Cells.Item(5, 40).Value = 3.9

# Row 6
$ws.Cells.Item(6, 20).Value = 1.39
$ws.Cells.Item(6, 21).Value = 2.38

# Row 7
$ws.Cells.Item(7, 7).Value = 1.21
$ws.Cells.Item(7, 9).Value = 22
$ws.Cells.Item(7, 10).Value = 8
$ws.Cells.Item(7, 11).Value = 9
$ws.Cells.Item(7, 19).Value = 2.4
$ws.Cells.Item(7, 20).Value = 2.48
$ws.Cells.Item(7, 22).Value = 1.05
$ws.Cells.Item(7, 23).Value = 5.6
$ws.Cells.Item(7, 24).Value = 34
$ws.Cells.Item(7, 26).Value = 1000
$ws.Cells.Item(7, 29).Value = 25
$ws.Cells.Item(7, 32).Value = 7.6
$ws.Cells.Item(7, 34).Value = 250
$ws.Cells.Item(7, 36).Value = 8.199999999999999
$ws.Cells.Item(7, 38).Value = 410
$ws.Cells.Item(7, 40).Value = 3.65

# Row 8
$ws.Cells.Item(8, 6).Value = 3.25
$ws.Cells.Item(8, 7).Value = 3.45
$ws.Cells.Item(8, 9).Value = 2.34
$ws.Cells.Item(8, 10).Value = 3.75
$ws.Cells.Item(8, 11).Value = 3.85

# Row 9
$ws.Cells.Item(9, 8).Value = 2.76
$ws.Cells.Item(9, 9).Value = 2.78
$ws.Cells.Item(9, 13).Value = 1.06
$ws.Cells.Item(9, 17).Value = 1.81
$ws.Cells.Item(9, 20).Value = 1.67
$ws.Cells.Item(9, 23).Value = 1.58
$ws.Cells.Item(9, 24).Value = 16.5
$ws.Cells.Item(9, 32).Value = 18.5
$ws.Cells.Item(9, 35).Value = 36
$ws.Cells.Item(9, 40).Value = 20

# Row 10
$ws.Cells.Item(10, 6).Value = 3.85
$ws.Cells.Item(10, 8).Value = 2
$ws.Cells.Item(10, 9).Value = 2.04
$ws.Cells.Item(10, 10).Value = 3.95
$ws.Cells.Item(10, 16).Value = 2.16
$ws.Cells.Item(10, 17).Value = 1.79
$ws.Cells.Item(10, 20).Value = 1.71
$ws.Cells.Item(10, 22).Value = 1.96
$ws.Cells.Item(10, 23).Value = 1.33
$ws.Cells.Item(10, 32).Value = 30
$ws.Cells.Item(10, 39).Value = 80

# Row 11
$ws.Cells.Item(11, 6).Value = 4.2
$ws.Cells.Item(11, 7).Value = 4.5
$ws.Cells.Item(11, 9).Value = 2.04
$ws.Cells.Item(11, 15).Value = 1.39
$ws.Cells.Item(11, 17).Value = 2.2
$ws.Cells.Item(11, 20).Value = 1.96
$ws.Cells.Item(11, 22).Value = 1.96
$ws.Cells.Item(11, 25).Value = 8.4

# Row 12
$ws.Cells.Item(12, 6).Value = 1.6
$ws.Cells.Item(12, 7).Value = 1.61
$ws.Cells.Item(12, 9).Value = 6.6
$ws.Cells.Item(12, 11).Value = 4.6
$ws.Cells.Item(12, 20).Value = 1.87
$ws.Cells.Item(12, 22).Value = 1.18
$ws.Cells.Item(12, 23).Value = 2.64
$ws.Cells.Item(12, 32).Value = 9.6
$ws.Cells.Item(12, 35).Value = 80
$ws.Cells.Item(12, 36).Value = 15

# Row 13
$ws.Cells.Item(13, 6).Value = 2.74
$ws.Cells.Item(13, 9).Value = 2.78
$ws.Cells.Item(13, 10).Value = 3.65
$ws.Cells.Item(13, 13).Value = 1.06
$ws.Cells.Item(13, 18).Value = 1.46
$ws.Cells.Item(13, 22).Value = 1.56
$ws.Cells.Item(13, 24).Value = 17.5
$ws.Cells.Item(13, 29).Value = 8
$ws.Cells.Item(13, 36).Value = 42

# Row 15
$ws.Cells.Item(15, 7).Value = 1.47
$ws.Cells.Item(15, 8).Value = 7.4
$ws.Cells.Item(15, 10).Value = 5.3
$ws.Cells.Item(15, 14).Value = 6.6
$ws.Cells.Item(15, 17).Value = 1.49
$ws.Cells.Item(15, 18).Value = 1.76
$ws.Cells.Item(15, 20).Value = 1.71
$ws.Cells.Item(15, 21).Value = 2.28
$ws.Cells.Item(15, 22).Value = 1.15
$ws.Cells.Item(15, 23).Value = 3.1
$ws.Cells.Item(15, 27).Value = 1000
$ws.Cells.Item(15, 28).Value = 13
$ws.Cells.Item(15, 29).Value = 13
$ws.Cells.Item(15, 30).Value = 30
$ws.Cells.Item(15, 32).Value = 11.5
$ws.Cells.Item(15, 33).Value = 11
$ws.Cells.Item(15, 36).Value = 14
$ws.Cells.Item(15, 37).Value = 14
$ws.Cells.Item(15, 38).Value = 26
$ws.Cells.Item(15, 39).Value = 80
$ws.Cells.Item(15, 40).Value = 4.8
$ws.Cells.Item(15, 41).Value = 80

# Row 16
$ws.Cells.Item(16, 7).Value = 4.5
$ws.Cells.Item(16, 8).Value = 2.28
$ws.Cells.Item(16, 9).Value = 2.48
$ws.Cells.Item(16, 14).Value = 2.2
$ws.Cells.Item(16, 15).Value = 1.71
$ws.Cells.Item(16, 20).Value = 2.42
$ws.Cells.Item(16, 21).Value = 1.6
$ws.Cells.Item(16, 22).Value = 1.68
$ws.Cells.Item(16, 35).Value = 110
$ws.Cells.Item(16, 36).Value = 150
$ws.Cells.Item(16, 37).Value = 110
$ws.Cells.Item(16, 38).Value = 160
$ws.Cells.Item(16, 39).Value = 380
$ws.Cells.Item(16, 40).Value = 200

# Row 17
$ws.Cells.Item(17, 6).Value = 2.46
$ws.Cells.Item(17, 12).Value = 1.3
$ws.Cells.Item(17, 15).Value = 1.21
$ws.Cells.Item(17, 16).Value = 2.48
$ws.Cells.Item(17, 18).Value = 1.6
$ws.Cells.Item(17, 19).Value = 2.54
$ws.Cells.Item(17, 20).Value = 1.56
$ws.Cells.Item(17, 24).Value = 22

# Row 18
$ws.Cells.Item(18, 12).Value = 1.4
$ws.Cells.Item(18, 14).Value = 3.9
$ws.Cells.Item(18, 19).Value = 3.5

# Row 19
$ws.Cells.Item(19, 12).Value = 1.41
$ws.Cells.Item(19, 14).Value = 3.8
$ws.Cells.Item(19, 18).Value = 1.35

# Row 20
$ws.Cells.Item(20, 7).Value = 2.1
$ws.Cells.Item(20, 8).Value = 3.8
$ws.Cells.Item(20, 12).Value = 1.37

# Row 21
$ws.Cells.Item(21, 6).Value = 1.9
$ws.Cells.Item(21, 7).Value = 1.95
$ws.Cells.Item(21, 12).Value = 1.32
$ws.Cells.Item(21, 14).Value = 5
$ws.Cells.Item(21, 17).Value = 1.69
$ws.Cells.Item(21, 19).Value = 2.74
$ws.Cells.Item(21, 27).Value = 85
$ws.Cells.Item(21, 29).Value = 9.4

# Row 22
$ws.Cells.Item(22, 7).Value = 2.92
$ws.Cells.Item(22, 9).Value = 2.68
$ws.Cells.Item(22, 12).Value = 1.38
$ws.Cells.Item(22, 14).Value = 4.1
$ws.Cells.Item(22, 17).Value = 1.89
$ws.Cells.Item(22, 19).Value = 3.25
$ws.Cells.Item(22, 22).Value = 1.59

# Row 23
$ws.Cells.Item(23, 12).Value = 1.41
$ws.Cells.Item(23, 14).Value = 4
$ws.Cells.Item(23, 17).Value = 1.97
$ws.Cells.Item(23, 19).Value = 3.5
$ws.Cells.Item(23, 21).Value = 2.26

# Row 24
$ws.Cells.Item(24, 7).Value = 2
$ws.Cells.Item(24, 8).Value = 4.1
$ws.Cells.Item(24, 9).Value = 4.3
$ws.Cells.Item(24, 10).Value = 3.8
$ws.Cells.Item(24, 12).Value = 1.37
$ws.Cells.Item(24, 14).Value = 4.2
$ws.Cells.Item(24, 17).Value = 1.85
$ws.Cells.Item(24, 19).Value = 3.15
$ws.Cells.Item(24, 23).Value = 2

# Row 26
$ws.Cells.Item(26, 6).Value = 2.38
$ws.Cells.Item(26, 12).Value = 1.4
$ws.Cells.Item(26, 14).Value = 3.95
$ws.Cells.Item(26, 16).Value = 1.98
$ws.Cells.Item(26, 17).Value = 1.93
$ws.Cells.Item(26, 19).Value = 3.35

# Row 27
$ws.Cells.Item(27, 6).Value = 2.58
$ws.Cells.Item(27, 9).Value = 3.65

# Row 28
$ws.Cells.Item(28, 8).Value = 4.5
